$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2390.6875
$ws.Range("I98").Value = 1207.3462
$ws.Range("J98").Value = 7518.5
$ws.Range("K98").Value = 1207.3462
$ws.Range("L98").Value = 7518.5
$ws.Range("M98").Value = 290.6538
$ws.Range("N98").Value = -10514.5
$ws.Range("H122").Value = 2390.6875
$ws.Range("I122").Value = 1207.3462
$ws.Range("J122").Value = 7518.5
$ws.Range("K122").Value = 3622.0386
$ws.Range("L122").Value = 22555.5
$ws.Range("M122").Value = -1172.0386
$ws.Range("N122").Value = -27455.5
$ws.Range("H138").Value = 3311.1355
$ws.Range("I138").Value = 1897.5834
$ws.Range("J138").Value = 3672.0425
$ws.Range("K138").Value = 5692.7502
$ws.Range("L138").Value = 11016.1275
$ws.Range("M138").Value = -552.7502000000004
$ws.Range("N138").Value = -21296.1275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5452.4097
$ws.Range("I32").Value = 3142.9285
$ws.Range("J32").Value = 10557.579
$ws.Range("K32").Value = 3142.9285
$ws.Range("L32").Value = 10557.579
$ws.Range("M32").Value = -2855.9285
$ws.Range("N32").Value = -11131.579
$ws.Range("H74").Value = 1608.2727
$ws.Range("I74").Value = 1265.6177
$ws.Range("J74").Value = 2773.3
$ws.Range("K74").Value = 1265.6177
$ws.Range("L74").Value = 2773.3
$ws.Range("M74").Value = -391.6177
$ws.Range("N74").Value = -4521.3
$ws.Range("H77").Value = 1608.2727
$ws.Range("I77").Value = 1265.6177
$ws.Range("J77").Value = 2773.3
$ws.Range("K77").Value = 6328.0885
$ws.Range("L77").Value = 13866.5
$ws.Range("M77").Value = -1960.0885
$ws.Range("N77").Value = -22602.5
$ws.Range("H137").Value = 41710
$ws.Range("J137").Value = 41710
$ws.Range("L137").Value = 41710
$ws.Range("N137").Value = -51910

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 118850
$ws.Range("J59").Value = 118850
$ws.Range("L59").Value = 118850
$ws.Range("N59").Value = -120544
$ws.Range("H87").Value = 41800
$ws.Range("J87").Value = 41800
$ws.Range("L87").Value = 41800
$ws.Range("N87").Value = -44296
$ws.Range("H90").Value = 41800
$ws.Range("J90").Value = 41800
$ws.Range("L90").Value = 125400
$ws.Range("N90").Value = -137880
$ws.Range("H137").Value = 35446.668
$ws.Range("J137").Value = 40670
$ws.Range("L137").Value = 40670
$ws.Range("N137").Value = -50870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27782144
$ws.Range("I31").Value = 2038.625
$ws.Range("J31").Value = 50006230
$ws.Range("K31").Value = 2038.625
$ws.Range("L31").Value = 50006230
$ws.Range("M31").Value = -1743.625
$ws.Range("N31").Value = -50006820
$ws.Range("H34").Value = 27782144
$ws.Range("I34").Value = 2038.625
$ws.Range("J34").Value = 50006230
$ws.Range("K34").Value = 2038.625
$ws.Range("L34").Value = 50006230
$ws.Range("M34").Value = -1836.625
$ws.Range("N34").Value = -50006634
$ws.Range("H99").Value = 9093998
$ws.Range("I99").Value = 22224512
$ws.Range("J99").Value = 3642.2307
$ws.Range("K99").Value = 22224512
$ws.Range("L99").Value = 3642.2307
$ws.Range("M99").Value = -22223014
$ws.Range("N99").Value = -6638.2307
$ws.Range("H126").Value = 9093998
$ws.Range("I126").Value = 22224512
$ws.Range("J126").Value = 3642.2307
$ws.Range("K126").Value = 66673536
$ws.Range("L126").Value = 10926.6921
$ws.Range("M126").Value = -66671066
$ws.Range("N126").Value = -15866.6921
$ws.Range("H132").Value = 3081.4285
$ws.Range("I132").Value = 2750.76
$ws.Range("J132").Value = 3908.1
$ws.Range("K132").Value = 8252.280000000001
$ws.Range("L132").Value = 11724.3
$ws.Range("M132").Value = -5722.280000000001
$ws.Range("N132").Value = -16784.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1850.5834
$ws.Range("I5").Value = 467.83334
$ws.Range("K5").Value = 1403.50002
$ws.Range("M5").Value = -1291.50002
$ws.Range("H23").Value = 180.86363
$ws.Range("I23").Value = 90.85714
$ws.Range("J23").Value = 222.86667
$ws.Range("K23").Value = 272.57142
$ws.Range("L23").Value = 668.60001
$ws.Range("M23").Value = -37.57141999999999
$ws.Range("N23").Value = -1138.60001
$ws.Range("H121").Value = 1799.5593
$ws.Range("J121").Value = 1799.5593
$ws.Range("L121").Value = 5398.6779
$ws.Range("N121").Value = -8018.6779
$ws.Range("H131").Value = 8929542
$ws.Range("J131").Value = 980.6415
$ws.Range("L131").Value = 2941.9245
$ws.Range("N131").Value = -13021.9245
$ws.Range("H132").Value = 1944.0682
$ws.Range("I132").Value = 790.7368
$ws.Range("J132").Value = 2820.6
$ws.Range("K132").Value = 7116.6312
$ws.Range("L132").Value = 25385.4
$ws.Range("M132").Value = -4586.6312
$ws.Range("N132").Value = -30445.4
$ws.Range("H135").Value = 1850.5834
$ws.Range("I135").Value = 467.83334
$ws.Range("K135").Value = 4210.50006
$ws.Range("M135").Value = -1675.50006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20836162
$ws.Range("I80").Value = 35716904
$ws.Range("K80").Value = 35716904
$ws.Range("M80").Value = -35715906
$ws.Range("H83").Value = 20836162
$ws.Range("I83").Value = 35716904
$ws.Range("K83").Value = 178584520
$ws.Range("M83").Value = -178579528
$ws.Range("H132").Value = 2899.516
$ws.Range("I132").Value = 1653.7059
$ws.Range("J132").Value = 4412.2856
$ws.Range("K132").Value = 4961.1177
$ws.Range("L132").Value = 13236.8568
$ws.Range("M132").Value = -2431.1177
$ws.Range("N132").Value = -18296.8568
$ws.Range("H137").Value = 45210
$ws.Range("J137").Value = 45210
$ws.Range("L137").Value = 45210
$ws.Range("N137").Value = -55410

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11116066
$ws.Range("I132").Value = 5869.8423
$ws.Range("J132").Value = 30306404
$ws.Range("K132").Value = 17609.5269
$ws.Range("L132").Value = 90919212
$ws.Range("M132").Value = -15079.5269
$ws.Range("N132").Value = -90924272
$ws.Range("H136").Value = 1320.3103
$ws.Range("I136").Value = 564.5
$ws.Range("J136").Value = 2999.889
$ws.Range("K136").Value = 1693.5
$ws.Range("L136").Value = 8999.667000000001
$ws.Range("M136").Value = 856.5
$ws.Range("N136").Value = -14099.667
